$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "WARNING: text not found for replacement: $old"
    }
}

# ---------------------------------------------------------------------------
# We work from the BOTTOM of the document upward so that paragraph indices
# used for the earlier (not-yet-processed) sections remain valid.
# ---------------------------------------------------------------------------

# --- Section: Risk Assessment ---------------------------------------------
# Paragraph 23: "No risk assessment is applicable at this time, as there are
# no pending updates available to assess."
# becomes one paragraph of body text, a blank paragraph, then three
# ListBullet paragraphs.
$p23 = $d.Paragraphs.Item(23)
$p23.Range.Text = "The current patches address some of the identified vulnerabilities; however, due to the nature of the updates and their availability, there is still a risk involved."

$p23.Range.InsertParagraphAfter()
$pBlank = $d.Paragraphs.Item(24)
$pBlank.Style = "Normal"

$pBlank.Range.InsertParagraphAfter()
$pRisk1 = $d.Paragraphs.Item(25)
$pRisk1.Range.Text = "**Potential Risk:** The potential risk is that the systems may not receive necessary security patches in a timely manner, leaving them vulnerable to attacks."
$pRisk1.Style = "ListBullet"

$pRisk1.Range.InsertParagraphAfter()
$pRisk2 = $d.Paragraphs.Item(26)
$pRisk2.Range.Text = "**Impact Level:** The impact level could be significant if the vulnerability were to be exploited, potentially causing system downtime or data breaches."
$pRisk2.Style = "ListBullet"

$pRisk2.Range.InsertParagraphAfter()
$pRisk3 = $d.Paragraphs.Item(27)
$pRisk3.Range.Text = "**Mitigation Plan:** To mitigate this risk, the organization should prioritize patching these updates as soon as possible."
$pRisk3.Style = "ListBullet"

# --- Section: Recommended next steps ---------------------------------------
# Paragraphs 18-20 (three numbered lines) collapse into a single paragraph.
$p19 = $d.Paragraphs.Item(19)
$p20 = $d.Paragraphs.Item(20)
$delRange = $d.Range($p19.Range.Start, $p20.Range.End)
$delRange.Delete()

$p18 = $d.Paragraphs.Item(18)
$p18.Range.Text = "The first step is to Review and Assess Updates to determine if they align with organizational requirements and identify potential risks. Next, Schedule patch deployments for the pending updates to ensure that the systems are brought up-to-date in a timely manner. Finally, update documentation to reflect any changes or additions made to system configurations."

# --- Section: Compliance with RMF Controls ---------------------------------
# Paragraph 11: simple text swap.
Replace-Text "In order to proceed with the RMF process for OS patch management, it's recommended that identification, reporting, and corrective action be conducted first. This involves checking for vulnerabilities and implementing patches as needed." "In order to meet the requirements of the Risk Management Framework (RMF), it is essential to identify and remediate vulnerabilities that could pose a risk to the system."

# Paragraph 12 stays blank. Paragraph 13 becomes a ListBullet item.
$p13 = $d.Paragraphs.Item(13)
$p13.Range.Text = "**Flaw Remediation:** The systems should be updated as soon as possible."
$p13.Style = "ListBullet"

# Paragraph 14 is an empty paragraph that must be removed entirely. After
# deletion "Vulnerability checks..." (originally paragraph 15) shifts down
# and becomes paragraph 14.
$p14 = $d.Paragraphs.Item(14)
$p14.Range.Delete()

# Paragraph 14 (now "Vulnerability checks...") turns into three ListBullet
# paragraphs.
$p14b = $d.Paragraphs.Item(14)
$p14b.Range.Text = "**Identification, Reporting / Corrective Action:** Identify potential risks and report them to management. Perform corrective actions, such as updating software or patching vulnerabilities, in a timely manner."
$p14b.Style = "ListBullet"

$p14b.Range.InsertParagraphAfter()
$pConfig = $d.Paragraphs.Item(15)
$pConfig.Range.Text = "**Configuration Management:** Ensure that configuration settings are properly documented and updated regularly."
$pConfig.Style = "ListBullet"

$pConfig.Range.InsertParagraphAfter()
$pVuln = $d.Paragraphs.Item(16)
$pVuln.Range.Text = "**Vulnerability Checks:** Regularly check the systems for known vulnerabilities and address any issues promptly."
$pVuln.Style = "ListBullet"

# --- Section: Patch Status Summary -----------------------------------------
# Paragraph 8 expands into an intro line, blank line, three ListBullet
# items, a blank line, and a closing line.
$p8 = $d.Paragraphs.Item(8)
$p8.Range.Text = "There are several systems with pending updates, including:"

$p8.Range.InsertParagraphAfter()
$pBlank2 = $d.Paragraphs.Item(9)
$pBlank2.Style = "Normal"

$pBlank2.Range.InsertParagraphAfter()
$pUpd1 = $d.Paragraphs.Item(10)
$pUpd1.Range.Text = " Linux-image-amd64: Updated from version 6.1.129-1 to 6.1.135-1."
$pUpd1.Style = "ListBullet"

$pUpd1.Range.InsertParagraphAfter()
$pUpd2 = $d.Paragraphs.Item(11)
$pUpd2.Range.Text = " perl-modules-5.36: Updated from version 5.36.0-7+deb12u1 to 5.36.0-7+deb12u2 (all)."
$pUpd2.Style = "ListBullet"

$pUpd2.Range.InsertParagraphAfter()
$pUpd3 = $d.Paragraphs.Item(12)
$pUpd3.Range.Text = " gir1.2-javascriptcoregtk-4.0, gir1.2-webkit2-4.0, libjavascriptcoregtk-4.0-18, and libwebkit2gtk-4.0: Updated from version 2.48.0-1~deb12u1 to 2.48.1-2~deb12u1 (all)."
$pUpd3.Style = "ListBullet"

$pUpd3.Range.InsertParagraphAfter()
$pBlank3 = $d.Paragraphs.Item(13)
$pBlank3.Style = "Normal"

$pBlank3.Range.InsertParagraphAfter()
$pClose = $d.Paragraphs.Item(14)
$pClose.Range.Text = "These updates are available for immediate installation."
$pClose.Style = "Normal"

# --- Section: System Overview -----------------------------------------------
Replace-Text "The system is currently undergoing updates to ensure compliance with existing regulations." "The systems in the network have received software updates to enhance security and stability."

# --- Created date ------------------------------------------------------------
Replace-Text "Created April 26, 2025 at 10:47:32" "Created April 27, 2025 at 15:55:13"

Write-Host "Done. Final paragraph count:" $d.Paragraphs.Count
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    Write-Host "$i : [$($p.Range.Text)]"
}
